$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet, row 3 (fa496268... entry):
#   Correspond Handoff Datetime (D3) and Correspond Handback DateTime (G3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-28 05:28:57"
$wsZhCn.Range("G3").Value = "2016-01-28 05:29:49"

# Update the "de-de" sheet, row 3 (fa496268... entry):
#   Correspond Handoff Datetime (D3) and Correspond Handback DateTime (G3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-28 05:29:08"
$wsDeDe.Range("G3").Value = "2016-01-28 05:30:10"
